# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets carry identical data in this workbook, so the same set of
# row -> new value updates is applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    7  = 468
    8  = 51
    10 = 587
    14 = 382
    16 = 101
    17 = 14
    19 = 53
    21 = 989
    22 = 1414
    24 = 338
    31 = 262
    33 = 1638
    39 = 3757
    40 = 436
    41 = 211
    42 = 928
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
